$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the "REGISTRE DIRECCIÓN" column (C) entirely; column D (NRO INFORME) shifts left to C.
$ws.Columns.Item(3).Delete()

# Rename the header to match the new meaning of the (shifted) column C.
$ws.Range("C1").Value = "NUMEROFICHA"

# Re-apply the (shifted) column width explicitly, since the delete can leave a split col run.
$ws.Columns.Item(3).ColumnWidth = 17.28515625

# Recompute the wrapped-text row heights now that the long address column is gone.
$ws.Rows.Item("1:23").AutoFit()

$ws.Range("G4").Select()
